# Update the two-digit division worksheet numbers.
# Each tuple is (old text, new text), listed in the same order the
# original values occur in the document. We search forward sequentially
# (without wrapping, without resetting the start position) so that
# duplicate values (e.g. "55÷6=" appears twice) are each replaced with
# their own distinct target text rather than a single blanket value.

$d = $word.ActiveDocument
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Forward = $true
$range.Find.Wrap = 0

$pairs = @(
    @("27÷9=", "29÷4="),
    @("14÷6=", "45÷5="),
    @("25÷3=", "83÷9="),
    @("49÷7=", "32÷8="),
    @("36÷3=", "74÷2="),
    @("84÷6=", "94÷8="),
    @("44÷7=", "11÷2="),
    @("89÷3=", "39÷9="),
    @("55÷6=", "84÷7="),
    @("24÷3=", "11÷5="),
    @("83÷4=", "55÷5="),
    @("83÷2=", "24÷6="),
    @("17÷8=", "77÷2="),
    @("77÷3=", "71÷8="),
    @("79÷8=", "64÷8="),
    @("29÷5=", "77÷7="),
    @("81÷5=", "75÷2="),
    @("50÷9=", "33÷4="),
    @("25÷8=", "49÷7="),
    @("26÷8=", "89÷4="),
    @("36÷8=", "62÷8="),
    @("57÷8=", "78÷8="),
    @("30÷5=", "45÷5="),
    @("55÷6=", "45÷2="),
    @("21÷5=", "38÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Text = $new
    } else {
        Write-Host "WARNING: could not find '$old'"
    }
}
